$wb = $excel.ActiveWorkbook

$guid = "a60b6ad0-a860-4fa0-9afc-39eb4409c937"
$commit = "52b2e1232c875aaef5436a650fd8309d0e3f146e"
$mdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$commit/e2e/$guid.md"

# ---------- Sheet "Overview" (sheet1, table3) ----------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A3").Value = "$guid.md"
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("D3").Borders.LineStyle = 0
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-18 22:39:36"
$wsOverview.Range("G3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOverview.Range("B3").Value = "e2e\$guid.md"
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $mdUrl, "", "", "e2e\$guid.md") | Out-Null

# ---------- Sheet "zh-cn" (sheet2, table1) ----------
$wsZh = $wb.Worksheets.Item("zh-cn")
$loZh = $wsZh.ListObjects.Item(1)
$loZh.ListRows.Add() | Out-Null

$wsZh.Range("B3").Value = ".md"
$wsZh.Range("C3").Value = "Ready for handoff"
$wsZh.Range("D3").Value = "e2e"
$wsZh.Range("E3").Value = "ht"
$wsZh.Range("F3").Value = "'False"
$wsZh.Range("G3").Value = "$guid.091ba2674379fee8428cbd33404325ac707c0e03.zh-cn.xlf"
$wsZh.Range("H3").Value = "2016-08-18 22:39:31"
$wsZh.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("I3").Borders.LineStyle = 0
$wsZh.Range("J3").Borders.LineStyle = 0
$wsZh.Range("K3").Value = "0001-01-01 00:00:00"
$wsZh.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("L3").Borders.LineStyle = 0
$wsZh.Range("M3").Value = "'True"
$wsZh.Range("N3").Borders.LineStyle = 0
$wsZh.Range("O3").Value = "'False"
$wsZh.Range("P3").Borders.LineStyle = 0

$wsZh.Range("A3").Value = "$guid.md"
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $mdUrl, "", "", "$guid.md") | Out-Null

# ---------- Sheet "de-de" (sheet3, table2) ----------
$wsDe = $wb.Worksheets.Item("de-de")
$loDe = $wsDe.ListObjects.Item(1)
$loDe.ListRows.Add() | Out-Null

$wsDe.Range("B3").Value = ".md"
$wsDe.Range("C3").Value = "Ready for handoff"
$wsDe.Range("D3").Value = "e2e"
$wsDe.Range("E3").Value = "ht"
$wsDe.Range("F3").Value = "'False"
$wsDe.Range("G3").Value = "$guid.091ba2674379fee8428cbd33404325ac707c0e03.de-de.xlf"
$wsDe.Range("H3").Value = "2016-08-18 22:39:36"
$wsDe.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("I3").Borders.LineStyle = 0
$wsDe.Range("J3").Borders.LineStyle = 0
$wsDe.Range("K3").Value = "0001-01-01 00:00:00"
$wsDe.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("L3").Borders.LineStyle = 0
$wsDe.Range("M3").Value = "'True"
$wsDe.Range("N3").Borders.LineStyle = 0
$wsDe.Range("O3").Value = "'False"
$wsDe.Range("P3").Borders.LineStyle = 0

$wsDe.Range("A3").Value = "$guid.md"
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $mdUrl, "", "", "$guid.md") | Out-Null

Write-Host "Edit complete"
